# Brewery Parts.xlsx edit script
# Commit: Renamed Lable Tun.lsl  Generalise Tun Contents for MT + Kettle

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Brewery")

# --- Cell value updates -----------------------------------------------

# Row 2 (Base)
$ws.Range("E2").Value2 = "Brewery"

# Row 4 (Mill)
$ws.Range("E4").Value2 = "Grain Mill"

# Row 5 (HLT)
$ws.Range("C5").Value2 = "Hot Liquor Tank"
$ws.Range("E5").Value2 = "Label Object from Desc"

# Row 6 (HLT Lid)
$ws.Range("E6").Value2 = "Tun Lid"
$ws.Range("I6").ClearContents()

# Row 7 (HLT contents)
$ws.Range("E7").Value2 = "HLT Steam"

# Row 9 (MT)
$ws.Range("E9").Value2 = "Label Object from Desc"

# Row 10 (MT Lid)
$ws.Range("E10").Value2 = "Tun Lid"
$ws.Range("I10").ClearContents()

# Row 11 (MT Contents)
$ws.Range("E11").Value2 = "Tun Contents"
$ws.Range("F11").Value2 = "Grain/Mash ??"
$ws.Range("I11").Value2 = "Differentiate grain to mash?"

# Row 13 (Kettle)
$ws.Range("E13").Value2 = "Label Object from Desc"

# Row 14 (Kettle lid)
$ws.Range("E14").Value2 = "Tun Lid"
$ws.Range("I14").ClearContents()

# Row 15 (Kettle contents)
$ws.Range("E15").Value2 = "Tun Contents"

# Row 16 (Chimney)
$ws.Range("I16").Value2 = "add puffer to top"

# Row 17 (Pipe1 / Pump Ass 1)
$ws.Range("E17").Value2 = "Pump Assembly 1"

# Row 18 (Pump1 / Pump Ass 1)
$ws.Range("E18").Value2 = "Pump Assembly 1`nLabel Obect from Desc"
$ws.Range("E18").WrapText = $true
$ws.Rows.Item(18).RowHeight = 31.5

# Row 19 (Pipe2a / Pump Ass 2)
$ws.Range("E19").Value2 = "Pump Assembly 2"

# Row 20 (pump2 / Pump Ass 2)
$ws.Range("E20").Value2 = "Pump Assembly 2`nLabel Object from Desc"
$ws.Range("E20").WrapText = $true
$ws.Rows.Item(20).RowHeight = 31.5
$ws.Range("I20").Interior.Color = 255

# Row 21 (pipe2b / Pump Ass 2)
$ws.Range("E21").Value2 = "Pump Assembly 2"
$ws.Range("I21").Interior.Color = 255

# Row 22 (pipe3a / Chiller circuit)
$ws.Range("E22").Value2 = "Chiller circuit"

# Row 23 (Chiller / Chiller circuit)
$ws.Range("E23").Value2 = "Chiller circuit`nLabel Object from Desc"
$ws.Range("E23").WrapText = $true
$ws.Rows.Item(23).RowHeight = 29.25

# Row 24 (pipe3b / Chiller circuit)
$ws.Range("E24").Value2 = "Chiller circuit"

# --- Column widths ------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 14.0
$ws.Columns.Item(5).ColumnWidth = 20.9

# --- Selection ------------------------------------------------------
$ws.Range("F11").Select() | Out-Null
